$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("A27").Value = 45860
$ws.Range("A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B27").Value = 687
$ws.Range("C27").Value = 781
$ws.Range("D27").Value = 687
$ws.Range("E27").Value = 781
$ws.Range("F27").Value = 2840
$ws.Range("G27").Value = 289

# Row 28
$ws.Range("A28").Value = 45861
$ws.Range("A28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B28").Value = 690
$ws.Range("C28").Value = 685.85
$ws.Range("D28").Value = 690
$ws.Range("E28").Value = 685.85
$ws.Range("F28").Value = 1499
$ws.Range("G28").Value = 275
